$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("Assigned") to make room for "Instance"
$ws.Columns.Item(3).Insert()

# Header row
$ws.Range("C1").Value = "Instance"

$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 'Paul, Matheo, Thor Waguespack'
$ws.Range("E2").Value = 'Paul, Matheo, Thor Waguespack'
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = 'Matheo, Blake Steel'
$ws.Range("E3").Value = 'Matheo, Blake Steel'
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 'Kamsi, Blake Steel, Noah Yaffe'
$ws.Range("E4").Value = 'Blake Steel, Noah Yaffe'
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = 'Jack Mogelof, Alejandro L, Jamari Pitchford'
$ws.Range("E5").Value = 'Jack Mogelof, Jamari Pitchford'
$ws.Range("C6").Value = 4
$ws.Range("D6").Value = 'Ali Awada, Alejandro E. Ulvert, Harry Corbin'
$ws.Range("E6").Value = 'Alejandro E. Ulvert, Harry Corbin'
$ws.Range("C7").Value = 5
$ws.Range("D7").Value = 'Kamsi, Jamari Pitchford'
$ws.Range("E7").Value = 'Jamari Pitchford'
$ws.Range("C8").Value = 6
$ws.Range("D8").Value = 'Jake Dieterich, Alejandro E. Ulvert'
$ws.Range("E8").Value = 'Alejandro E. Ulvert'
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 'Kamsi, Edu, Adi'
$ws.Range("E9").Value = 'Edu, Adi'
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 'Matheo, Noah Yaffe'
$ws.Range("E10").Value = 'Matheo, Noah Yaffe'
$ws.Range("C11").Value = 2
$ws.Range("D11").Value = 'Alexander, Alejandro L, George Ryckman'
$ws.Range("E11").Value = 'Alexander, George Ryckman'
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 'Jaxon, Kamsi, George Ryckman'
$ws.Range("E12").Value = 'Jaxon, George Ryckman'
$ws.Range("C13").Value = 4
$ws.Range("D13").Value = 'Jamari Pitchford, Edu'
$ws.Range("E13").Value = 'Jamari Pitchford, Edu'
$ws.Range("C14").Value = 5
$ws.Range("D14").Value = 'Matheo, Alejandro Espinosa'
$ws.Range("E14").Value = 'Matheo, Alejandro Espinosa'
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 'Josh Greene, Alejandro Espinosa, Jamari Pitchford'
$ws.Range("E15").Value = 'Josh Greene, Alejandro Espinosa, Jamari Pitchford'
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 'Thor Waguespack, Ali Awada'
$ws.Range("E16").Value = 'Thor Waguespack'
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 'Jack Mogelof, Alejandro L, Adi'
$ws.Range("E17").Value = 'Jack Mogelof, Adi'
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 'Gabe Heller, Alejandro L, Adi'
$ws.Range("E18").Value = 'Gabe Heller, Adi'
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 'Gabe Heller, Jake Dieterich, Blake Steel'
$ws.Range("E19").Value = 'Gabe Heller, Blake Steel'
$ws.Range("C20").Value = 5
$ws.Range("D20").Value = 'Paul, Henry, Harry Corbin'
$ws.Range("E20").Value = 'Paul, Harry Corbin'
$ws.Range("C21").Value = 6
$ws.Range("D21").Value = 'Ezana, Edu'
$ws.Range("E21").Value = 'Edu'
$ws.Range("C22").Value = 0
$ws.Range("D22").Value = 'Paul, Josh Greene, George Ryckman'
$ws.Range("E22").Value = 'Paul, Josh Greene, George Ryckman'
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 'Alejandro Espinosa, Thor Waguespack'
$ws.Range("E23").Value = 'Alejandro Espinosa, Thor Waguespack'
$ws.Range("C24").Value = 2
$ws.Range("D24").Value = 'Henry, Adi, Noah Yaffe'
$ws.Range("E24").Value = 'Adi, Noah Yaffe'
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 'Ben Kairouz, Alejandro Espinosa, Alejandro E. Ulvert'
$ws.Range("E25").Value = 'Ben Kairouz, Alejandro Espinosa, Alejandro E. Ulvert'
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 'Gabe Heller, Alexander, Ali Awada'
$ws.Range("E26").Value = 'Gabe Heller, Alexander'
$ws.Range("C27").Value = 5
$ws.Range("D27").Value = 'Jack Mogelof, Jake Dieterich'
$ws.Range("E27").Value = 'Jack Mogelof'
$ws.Range("C28").Value = 0
$ws.Range("D28").Value = 'Thor Waguespack, Alejandro E. Ulvert, Noah Yaffe'
$ws.Range("E28").Value = 'Thor Waguespack, Alejandro E. Ulvert, Noah Yaffe'
$ws.Range("C29").Value = 1
$ws.Range("D29").Value = 'Jack Mogelof, Ali Awada'
$ws.Range("E29").Value = 'Jack Mogelof'
$ws.Range("C30").Value = 2
$ws.Range("D30").Value = 'Josh Greene, Henry, Adi'
$ws.Range("E30").Value = 'Josh Greene, Adi'
$ws.Range("C31").Value = 3
$ws.Range("D31").Value = 'Thor Waguespack, Blake Steel, Harry Corbin'
$ws.Range("E31").Value = 'Thor Waguespack, Blake Steel, Harry Corbin'
$ws.Range("C32").Value = 4
$ws.Range("D32").Value = 'Josh Greene, Henry, George Ryckman'
$ws.Range("E32").Value = 'Josh Greene, George Ryckman'
$ws.Range("C33").Value = 5
$ws.Range("D33").Value = 'Alexander, Edu'
$ws.Range("E33").Value = 'Alexander, Edu'
$ws.Range("C34").Value = 6
$ws.Range("D34").Value = 'Jaxon, Alejandro Espinosa'
$ws.Range("E34").Value = 'Jaxon, Alejandro Espinosa'
$ws.Range("C35").Value = 0
$ws.Range("D35").Value = 'Gabe Heller, Ezana, Harry Corbin'
$ws.Range("E35").Value = 'Gabe Heller, Harry Corbin'
$ws.Range("C36").Value = 1
$ws.Range("D36").Value = 'Jaxon, Ezana, Ali Awada'
$ws.Range("E36").Value = 'Jaxon'
$ws.Range("C37").Value = 2
$ws.Range("D37").Value = 'Jaxon, Gabe Heller, Alexander'
$ws.Range("E37").Value = 'Jaxon, Gabe Heller, Alexander'
$ws.Range("C38").Value = 3
$ws.Range("D38").Value = 'Ben Kairouz, Alexander, Harry Corbin'
$ws.Range("E38").Value = 'Ben Kairouz, Alexander, Harry Corbin'
$ws.Range("C39").Value = 4
$ws.Range("D39").Value = 'Alejandro L, George Ryckman, Blake Steel'
$ws.Range("E39").Value = 'George Ryckman, Blake Steel'
$ws.Range("C40").Value = 5
$ws.Range("D40").Value = 'Ezana, Edu'
$ws.Range("E40").Value = 'Edu'
$ws.Range("C41").Value = 6
$ws.Range("D41").Value = 'Henry, Jamari Pitchford'
$ws.Range("E41").Value = 'Jamari Pitchford'
$ws.Range("C42").Value = 0
$ws.Range("D42").Value = 'Ben Kairouz, Kamsi, Paul'
$ws.Range("E42").Value = 'Ben Kairouz, Paul'
$ws.Range("C43").Value = 1
$ws.Range("D43").Value = 'Jaxon, Jake Dieterich'
$ws.Range("E43").Value = 'Jaxon'
$ws.Range("C44").Value = 2
$ws.Range("D44").Value = 'Jack Mogelof, Josh Greene, Jake Dieterich'
$ws.Range("E44").Value = 'Jack Mogelof, Josh Greene'
$ws.Range("C45").Value = 3
$ws.Range("D45").Value = 'Paul, Matheo, Noah Yaffe'
$ws.Range("E45").Value = 'Paul, Matheo, Noah Yaffe'
$ws.Range("C46").Value = 4
$ws.Range("D46").Value = 'Ben Kairouz, Ezana'
$ws.Range("E46").Value = 'Ben Kairouz'
$ws.Range("C47").Value = 5
$ws.Range("D47").Value = 'Ben Kairouz, Alejandro E. Ulvert'
$ws.Range("E47").Value = 'Ben Kairouz, Alejandro E. Ulvert'
